# Scheduled-runner update: refresh computed market/profit columns
# (currentAveragePrice / LevePrice / LeveProfit, columns H:N) on each
# class sheet, matching freshly pulled Universalis price data.
#
# Columns per sheet:
#   H currentAveragePrice   I currentAveragePriceNQ   J currentAveragePriceHQ
#   K LevePriceNQ           L LevePriceHQ
#   M LeveProfitNQ          N LeveProfitHQ
# Rows only carry the profit column(s) relevant to whether the leve's
# reward is NQ-priced, HQ-priced, or both, so a refresh can add/drop the
# M or N cell for a row depending on which side came back nonzero.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 573.5
$ws.Range("I101").Value = 573.5
$ws.Range("K101").Value = 1720.5
$ws.Range("M101").Value = -98.5

$ws.Range("H106").Value = 2735.2942
$ws.Range("I106").Value = 1983.7778
$ws.Range("J106").Value = 3580.75
$ws.Range("K106").Value = 1983.7778
$ws.Range("L106").Value = 3580.75
$ws.Range("M106").Value = -1352.7778
$ws.Range("N106").Value = -4842.75

$ws.Range("H107").Value = 947.2273
$ws.Range("I107").Value = 1011.4375
$ws.Range("J107").Value = 776
$ws.Range("K107").Value = 1011.4375
$ws.Range("L107").Value = 776
$ws.Range("M107").Value = 908.5625
$ws.Range("N107").Value = -4616

$ws.Range("H129").Value = 824.40985
$ws.Range("J129").Value = 849.3090999999999
$ws.Range("L129").Value = 2547.9273
$ws.Range("N129").Value = -12547.9273

$ws.Range("H132").Value = 3594.1304
$ws.Range("I132").Value = 3530.2273
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 10590.6819
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -8060.6819
$ws.Range("N132").Value = -20060

$ws.Range("H137").Value = 149623.8
$ws.Range("I137").Value = 238336.36
$ws.Range("J137").Value = 60911.234
$ws.Range("K137").Value = 715009.08
$ws.Range("L137").Value = 182733.702
$ws.Range("M137").Value = -712459.08
$ws.Range("N137").Value = -187833.702

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22272.191
$ws.Range("I32").Value = 25677.697
$ws.Range("J32").Value = 6001.4443
$ws.Range("K32").Value = 25677.697
$ws.Range("L32").Value = 6001.4443
$ws.Range("M32").Value = -25390.697
$ws.Range("N32").Value = -6575.4443

$ws.Range("H61").Value = 5332.6665
$ws.Range("I61").Value = 2999.5
$ws.Range("J61").Value = 5999.2856
$ws.Range("K61").Value = 2999.5
$ws.Range("L61").Value = 5999.2856
$ws.Range("M61").Value = -2787.5
$ws.Range("N61").Value = -6423.2856

$ws.Range("H63").Value = 2405802.5
$ws.Range("J63").Value = 5209955
$ws.Range("L63").Value = 5209955
$ws.Range("N63").Value = -5211327

$ws.Range("H66").Value = 2405802.5
$ws.Range("J66").Value = 5209955
$ws.Range("L66").Value = 26049775
$ws.Range("N66").Value = -26056639

$ws.Range("H74").Value = 2270.353
$ws.Range("I74").Value = 1748.5
$ws.Range("K74").Value = 1748.5
$ws.Range("M74").Value = -874.5

$ws.Range("H77").Value = 2270.353
$ws.Range("I77").Value = 1748.5
$ws.Range("K77").Value = 8742.5
$ws.Range("M77").Value = -4374.5

$ws.Range("H132").Value = 18292.549
$ws.Range("J132").Value = 48746.547
$ws.Range("L132").Value = 146239.641
$ws.Range("N132").Value = -151299.641

$ws.Range("H136").Value = 5332.6665
$ws.Range("I136").Value = 2999.5
$ws.Range("J136").Value = 5999.2856
$ws.Range("K136").Value = 8998.5
$ws.Range("L136").Value = 17997.8568
$ws.Range("M136").Value = -6448.5
$ws.Range("N136").Value = -23097.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 434.5
$ws.Range("I22").Value = 461.77777
$ws.Range("K22").Value = 461.77777
$ws.Range("M22").Value = -288.77777

$ws.Range("H86").Value = 1837.4166
$ws.Range("I86").Value = 1661.0555
$ws.Range("J86").Value = 2366.5
$ws.Range("K86").Value = 1661.0555
$ws.Range("L86").Value = 2366.5
$ws.Range("M86").Value = -538.0554999999999
$ws.Range("N86").Value = -4612.5

$ws.Range("H89").Value = 1837.4166
$ws.Range("I89").Value = 1661.0555
$ws.Range("J89").Value = 2366.5
$ws.Range("K89").Value = 8305.2775
$ws.Range("L89").Value = 11832.5
$ws.Range("M89").Value = -2689.2775
$ws.Range("N89").Value = -23064.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11849.795
$ws.Range("I31").Value = 21303.938
$ws.Range("J31").Value = 5273
$ws.Range("K31").Value = 21303.938
$ws.Range("L31").Value = 5273
$ws.Range("M31").Value = -21008.938
$ws.Range("N31").Value = -5863

$ws.Range("H34").Value = 11849.795
$ws.Range("I34").Value = 21303.938
$ws.Range("J34").Value = 5273
$ws.Range("K34").Value = 21303.938
$ws.Range("L34").Value = 5273
$ws.Range("M34").Value = -21101.938
$ws.Range("N34").Value = -5677

$ws.Range("H132").Value = 21934.691
$ws.Range("I132").Value = 24531.318
$ws.Range("J132").Value = 7653.25
$ws.Range("K132").Value = 73593.954
$ws.Range("L132").Value = 22959.75
$ws.Range("M132").Value = -71063.954
$ws.Range("N132").Value = -28019.75

$ws.Range("H134").Value = 1275.0769
$ws.Range("I134").Value = 975.1111
$ws.Range("J134").Value = 1950
$ws.Range("K134").Value = 2925.3333
$ws.Range("L134").Value = 5850
$ws.Range("M134").Value = -390.3332999999998
$ws.Range("N134").Value = -10920

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 725.9091
$ws.Range("J92").Value = 869.2857
$ws.Range("L92").Value = 2607.8571
$ws.Range("N92").Value = -5103.8571

$ws.Range("H131").Value = 121320.695
$ws.Range("I131").Value = 1010
$ws.Range("J131").Value = 125832.35
$ws.Range("K131").Value = 3030
$ws.Range("L131").Value = 377497.05
$ws.Range("M131").Value = 2010
$ws.Range("N131").Value = -387577.05

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 3000
$ws.Range("I6").Value = 3000
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 3000
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -2887
$ws.Range("N6").ClearContents()

$ws.Range("H16").Value = 3000
$ws.Range("I16").Value = 3000
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 3000
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -2750
$ws.Range("N16").ClearContents()

$ws.Range("H132").Value = 87267.39
$ws.Range("I132").Value = 104482.5
$ws.Range("J132").Value = 65748.5
$ws.Range("K132").Value = 313447.5
$ws.Range("L132").Value = 197245.5
$ws.Range("M132").Value = -310917.5
$ws.Range("N132").Value = -202305.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 982.381
$ws.Range("I46").Value = 714.1177
$ws.Range("J46").Value = 2122.5
$ws.Range("K46").Value = 714.1177
$ws.Range("L46").Value = 2122.5
$ws.Range("M46").Value = -526.1177
$ws.Range("N46").Value = -2498.5

$ws.Range("H93").Value = 2048.6
$ws.Range("I93").Value = 2279.1177
$ws.Range("J93").Value = 1558.75
$ws.Range("K93").Value = 2279.1177
$ws.Range("L93").Value = 1558.75
$ws.Range("M93").Value = -1031.1177
$ws.Range("N93").Value = -4054.75

$ws.Range("H136").Value = 25636.318
$ws.Range("I136").Value = 35472.266
$ws.Range("K136").Value = 106416.798
$ws.Range("M136").Value = -103866.798

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 8333.333000000001
$ws.Range("I48").Value = 8333.333000000001
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 8333.333000000001
$ws.Range("L48").Value = 0
$ws.Range("M48").Value = -7764.333000000001
$ws.Range("N48").ClearContents()

$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

$ws.Range("H100").Value = 1328.3
$ws.Range("I100").Value = 770
$ws.Range("J100").Value = 1700.5
$ws.Range("K100").Value = 1540
$ws.Range("L100").Value = 3401
$ws.Range("M100").Value = -999
$ws.Range("N100").Value = -4483
